# "Generate Report for Archive"
#
# The localization-status report is regenerated:
#   - every "Ready for handoff" status cell becomes "In Translation"
#   - the Status-related columns are narrowed to match the refreshed report layout
#     (Overview!E:F, and column C on each per-language detail sheet)

$wb = $excel.ActiveWorkbook

# --- Update status text everywhere it appears (Overview E/F, and the "Status"
#     column on each language sheet) without having to hardcode every address ---
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- Narrow the status columns to the refreshed report width ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5       # column C (Status)
